$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2..515) holds the "Förändrad" (changed) date.
# Update the serial date value from 45188 (2023-09-19) to 45189 (2023-09-20)
# for every data row, leaving everything else untouched.
$ws.Range("C2:C515").Value = 45189
